$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 9700
$ws.Range("J7").Value = 9700
$ws.Range("L7").Value = 9700
$ws.Range("N7").Value = -9924

$ws.Range("H14").Value = 9700
$ws.Range("J14").Value = 9700
$ws.Range("L14").Value = 9700
$ws.Range("N14").Value = -10082

$ws.Range("H17").Value = 383236
$ws.Range("J17").Value = 424645.4
$ws.Range("L17").Value = 1273936.2
$ws.Range("N17").Value = -1274272.2

$ws.Range("H40").Value = 1303.5555
$ws.Range("I40").Value = 1188
$ws.Range("J40").Value = 1448
$ws.Range("K40").Value = 1188
$ws.Range("L40").Value = 1448
$ws.Range("M40").Value = -1013
$ws.Range("N40").Value = -1798

$ws.Range("H51").Value = 2438.4285
$ws.Range("I51").Value = 2140.6
$ws.Range("J51").Value = 2488.0667
$ws.Range("K51").Value = 2140.6
$ws.Range("L51").Value = 2488.0667
$ws.Range("M51").Value = -1656.6
$ws.Range("N51").Value = -3456.0667

$ws.Range("H112").Value = 2364839.5
$ws.Range("J112").Value = 3468031.2
$ws.Range("L112").Value = 10404093.6
$ws.Range("N112").Value = -10406309.6

$ws.Range("H125").Value = 6662.6665
$ws.Range("I125").Value = 24988.8
$ws.Range("J125").Value = 935.75
$ws.Range("K125").Value = 224899.2
$ws.Range("L125").Value = 8421.75
$ws.Range("M125").Value = -222439.2
$ws.Range("N125").Value = -13341.75

$ws.Range("H129").Value = 688.7692
$ws.Range("I129").Value = 535
$ws.Range("J129").Value = 934.8
$ws.Range("K129").Value = 1605
$ws.Range("L129").Value = 2804.4
$ws.Range("M129").Value = 3395
$ws.Range("N129").Value = -12804.4

$ws.Range("H132").Value = 2298.5874
$ws.Range("I132").Value = 1894.5438
$ws.Range("K132").Value = 5683.6314
$ws.Range("M132").Value = -3153.6314

$ws.Range("H135").Value = 346.90475
$ws.Range("I135").Value = 312.125
$ws.Range("K135").Value = 2809.125
$ws.Range("M135").Value = -274.125

$ws.Range("H138").Value = 2377.797
$ws.Range("J138").Value = 3896.5483
$ws.Range("L138").Value = 11689.6449
$ws.Range("N138").Value = -21969.6449

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3918.19
$ws.Range("I32").Value = 3008.8877
$ws.Range("J32").Value = 11275.272
$ws.Range("K32").Value = 3008.8877
$ws.Range("L32").Value = 11275.272
$ws.Range("M32").Value = -2721.8877
$ws.Range("N32").Value = -11849.272

$ws.Range("H37").Value = 11044.667
$ws.Range("I37").Value = 7122.6665
$ws.Range("J37").Value = 14966.667
$ws.Range("K37").Value = 7122.6665
$ws.Range("L37").Value = 14966.667
$ws.Range("M37").Value = -6849.6665
$ws.Range("N37").Value = -15512.667

$ws.Range("H44").Value = 24450
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 24450
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 24450
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -25426

$ws.Range("H55").Value = 25000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 25000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 25000
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -25630

$ws.Range("H61").Value = 1960.1818
$ws.Range("I61").Value = 1880.4166
$ws.Range("J61").Value = 2055.9
$ws.Range("K61").Value = 1880.4166
$ws.Range("L61").Value = 2055.9
$ws.Range("M61").Value = -1668.4166
$ws.Range("N61").Value = -2479.9

$ws.Range("H136").Value = 1960.1818
$ws.Range("I136").Value = 1880.4166
$ws.Range("J136").Value = 2055.9
$ws.Range("K136").Value = 5641.2498
$ws.Range("L136").Value = 6167.700000000001
$ws.Range("M136").Value = -3091.2498
$ws.Range("N136").Value = -11267.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1154.4736
$ws.Range("I134").Value = 961.4074000000001
$ws.Range("J134").Value = 1628.3636
$ws.Range("K134").Value = 2884.2222
$ws.Range("L134").Value = 4885.0908
$ws.Range("M134").Value = -349.2222000000002
$ws.Range("N134").Value = -9955.0908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 5358.3335
$ws.Range("I26").Value = 1000
$ws.Range("J26").Value = 7537.5
$ws.Range("K26").Value = 1000
$ws.Range("L26").Value = 7537.5
$ws.Range("M26").Value = -713
$ws.Range("N26").Value = -8111.5

$ws.Range("H31").Value = 36102.195
$ws.Range("I31").Value = 1139.5
$ws.Range("K31").Value = 1139.5
$ws.Range("M31").Value = -844.5

$ws.Range("H34").Value = 36102.195
$ws.Range("I34").Value = 1139.5
$ws.Range("K34").Value = 1139.5
$ws.Range("M34").Value = -937.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4791.393
$ws.Range("I56").Value = 4791.393
$ws.Range("K56").Value = 4791.393
$ws.Range("M56").Value = -4261.393

$ws.Range("H131").Value = 842.0923
$ws.Range("J131").Value = 990.5238000000001
$ws.Range("L131").Value = 2971.5714
$ws.Range("N131").Value = -13051.5714

$ws.Range("H133").Value = 375282.22
$ws.Range("J133").Value = 593842.9399999999
$ws.Range("L133").Value = 1781528.82
$ws.Range("N133").Value = -1791648.82

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4910.0835
$ws.Range("I132").Value = 5623.643
$ws.Range("J132").Value = 2412.625
$ws.Range("K132").Value = 16870.929
$ws.Range("L132").Value = 7237.875
$ws.Range("M132").Value = -14340.929
$ws.Range("N132").Value = -12297.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10969.4
$ws.Range("I22").Value = 900
$ws.Range("J22").Value = 17682.334
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 17682.334
$ws.Range("M22").Value = -605
$ws.Range("N22").Value = -18272.334

$ws.Range("H27").Value = 10969.4
$ws.Range("I27").Value = 900
$ws.Range("J27").Value = 17682.334
$ws.Range("K27").Value = 900
$ws.Range("L27").Value = 17682.334
$ws.Range("M27").Value = -793
$ws.Range("N27").Value = -17896.334

$ws.Range("H46").Value = 2575.1667
$ws.Range("I46").Value = 10001
$ws.Range("J46").Value = 1090
$ws.Range("K46").Value = 10001
$ws.Range("L46").Value = 1090
$ws.Range("M46").Value = -9813
$ws.Range("N46").Value = -1466

$ws.Range("H55").Value = 52631810
$ws.Range("I55").Value = 83333460
$ws.Range("J55").Value = 415.57144
$ws.Range("K55").Value = 83333460
$ws.Range("L55").Value = 415.57144
$ws.Range("M55").Value = -83333287
$ws.Range("N55").Value = -761.5714399999999

$ws.Range("H136").Value = 3756.625
$ws.Range("I136").Value = 1503.325
$ws.Range("J136").Value = 15023.125
$ws.Range("K136").Value = 4509.975
$ws.Range("L136").Value = 45069.375
$ws.Range("M136").Value = -1959.975
$ws.Range("N136").Value = -50169.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 6000
$ws.Range("J52").Value = 6000
$ws.Range("L52").Value = 6000
$ws.Range("N52").Value = -6452

$ws.Range("H132").Value = 3010.423
$ws.Range("I132").Value = 3660.6943
$ws.Range("J132").Value = 1547.3125
$ws.Range("K132").Value = 10982.0829
$ws.Range("L132").Value = 4641.9375
$ws.Range("M132").Value = -8452.082900000001
$ws.Range("N132").Value = -9701.9375
